$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# This reflects a reshuffle of the weekly price records (row 5 unchanged).

$rows = @(
    @{ Row = 2;  D = 44466; J = 400; K = 9500;  L = 10000; M = 9750;  P = 390 },
    @{ Row = 3;  D = 44377; J = 520; K = 12500; L = 13000; M = 12750; P = 510 },
    @{ Row = 4;  D = 44372; J = 500; K = 13000; L = 14000; M = 13500; P = 540 },
    @{ Row = 6;  D = 44484; J = 400; K = 9000;  L = 10000; M = 9500;  P = 380 },
    @{ Row = 7;  D = 44384; J = 560; K = 11500; L = 12000; M = 11750; P = 470 },
    @{ Row = 8;  D = 44376; J = 400; K = 12000; L = 13000; M = 12500; P = 500 },
    @{ Row = 9;  D = 44316; J = 300; K = 16000; L = 17000; M = 16500; P = 660 },
    @{ Row = 10; D = 44370; J = 520; K = 13000; L = 14000; M = 13500; P = 540 },
    @{ Row = 11; D = 44425; J = 400; K = 11500; L = 12000; M = 11750; P = 470 },
    @{ Row = 12; D = 44386; J = 500; K = 11000; L = 12000; M = 11500; P = 460 },
    @{ Row = 13; D = 44446; J = 500; K = 11000; L = 12000; M = 11500; P = 460 },
    @{ Row = 14; D = 44473; J = 500; K = 8500;  L = 9000;  M = 8750;  P = 350 }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("J$r").Value = $rec.J
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("P$r").Value = $rec.P
}
